$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in attendance marks ("p") for two additional days (columns O and P)
# across every student row (3-30), mirroring the existing columns E..N.
for ($row = 3; $row -le 30; $row++) {
    $ws.Cells.Item($row, 15).Value = "p"
    $ws.Cells.Item($row, 16).Value = "p"
}

# The two previously-visible "day" columns immediately before the new data
# (L and M) are now hidden, same as the other already-hidden day columns.
$ws.Columns("L:M").Hidden = $true

# Row 9 had picked up a slightly-off custom height; restore it to the
# sheet's standard row height now that the row has been touched again.
$ws.Rows.Item(9).UseStandardHeight = $true
$ws.Rows.Item(9).AutoFit()

# Reflect where the user ended up after entering the new data.
$ws.Range("P33").Select() | Out-Null
